# Fruta / hortaliza, semanal
#
# A new weekly price row is inserted into the "Mandarina" sheet at row 140,
# pushing the existing rows 140-161 down to 141-162 (all of their values,
# styles, and number formats move with them). The freshly inserted row 140
# is then populated with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 140; Excel shifts rows 140:161
# down to 141:162 (values + formatting included).
$ws.Rows.Item(140).Insert()

# Populate the newly inserted row 140 with the new weekly record.
$ws.Range("A140").Value = 7
$ws.Range("B140").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C140").Value = "Ñuble"
$ws.Range("D140").Value = 44505
$ws.Range("E140").Value = 16
$ws.Range("F140").Value = "Fruta"
$ws.Range("G140").Value = 100102
$ws.Range("H140").Value = "Cítricos"
$ws.Range("I140").Value = 100102004
$ws.Range("J140").Value = "Mandarina"
$ws.Range("K140").Value = "Murcott"
$ws.Range("L140").Value = "Primera"
$ws.Range("M140").Value = 120
$ws.Range("N140").Value = 6000
$ws.Range("O140").Value = 6500
$ws.Range("P140").Value = 6250
$ws.Range("Q140").Value = "`$/caja 18 kilos"
$ws.Range("R140").Value = "Región de O'Higgins"
$ws.Range("S140").Value = 347
$ws.Range("T140").Value = 18

# Note: Rows.Insert() already carries the surrounding row's formatting
# (including column D's date/time number format) onto the new row 140,
# so no extra formatting step is required here.
